$d = $word.ActiveDocument

function Replace-One([string]$OldText, [string]$NewText, [object]$BoldVal, [object]$ColorVal, [object]$UnderlineVal) {
    $r = $d.Content
    $found = $r.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 1)
    if (-not $found) {
        Write-Host "NOT FOUND:" $OldText
        return
    }
    if ($NewText -eq "") {
        return
    }
    if ($BoldVal -ne $null) {
        $r.Font.Bold = $BoldVal
    }
    if ($ColorVal -ne $null) {
        $r.Font.Color = $ColorVal
    }
    if ($UnderlineVal -ne $null) {
        $r.Font.Underline = $UnderlineVal
    }
}

Replace-One "ENGLISH / " "英语/ " $null $null $null
Replace-One "SPANISH" "西班牙语" $null 13391121 1
Replace-One "PORTUGUESE " "葡萄牙语" $null 13391121 1
Replace-One "FRENCH " "法语" $null 13391121 1
Replace-One "VIETNAMESE " "越南语" $null 13391121 1
Replace-One "THAI" "泰语" $null 13391121 1
Replace-One "Background:" "背景：" 1 $null $null
Replace-One "This is an invitation for affiliates to a meeting in their country or city. This is the first email to go out " "这是一份邀请函，邀请联盟会员参加在其所在国家或城市举行的会议。 这是第一封发出的电子邮件 " $null $null $null
Replace-One "ENGLISH" "英语" $null $null $null
Replace-One "Subject line:" "主题：" 1 $null $null
Replace-One "Deriv Affiliate meeting | [CITY NAME] | [DATE]" "Deriv 联盟会议 | [城市名] | [DATE]" $null $null $null
Replace-One "Body:" "正文:" 1 $null $null
Replace-One "See you in [CITY NAME]!" "在 [城市名] 见！" $null $null $null
Replace-One "Great news! We will be in [CITY NAME] from" "好消息！ 我们将于" $null $null $null
Replace-One "[DATE] to [DATE] 2023. Our affiliate team, led by [COUNTRY] Country Manager [AFFILIATE MANAGER NAME], look forward to an exclusive one-on-one session with you." " 2023 年 [DATE] 到 [DATE] 抵达[城市名]。 我们的联盟团队由 [COUNTRY] 区域经理 [联盟经理名称] 领导，期待与您一对一会话。" $null $null $null
Replace-One "We’d love to hear about your experience with our affiliate programme. If there’s any way we can improve your experience, here’s your chance to tell us." "我们很想听听您参与联盟计划的经历。 如果有什么方法可以改善您的体验，请乘此机会告诉我们。" $null $null $null
Replace-One "When?" "什么时候？" $null $null $null
Replace-One "A 1-hour slot between 9:00 AM and 6:00 PM" "" $null $null $null
Replace-One "from [DATE] to [DATE]" "从 [DATE] 到 [DATE] 上午 9:00 至下午 6:00 之间 1 小时的时段" $null $null $null
Replace-One "Where?" "在哪里？" 1 $null $null
Replace-One "To be confirmed" "待确认" $null $null $null
Replace-One "How to book a slot?" "如何预订时段？" $null $null $null
Replace-One "Pick a date and time, and reply to this email by [DATE]  (first come, first served)" "请选择日期和时间，并在 [DATE]  前回复此电子邮件（先到先得）" $null $null $null
Replace-One "You’re welcome to bring along your clients and friends interested in learning more about trading on Deriv." "欢迎带上客户和有兴趣更加了解 Deriv 交易的朋友。" $null $null $null
Replace-One "We’re grateful for your continuous support and look forward to meeting you!" "非常感谢您一直以来的支持，并期待与您见面！" $null $null $null
Replace-One "P.S. We’re giving out free Deriv merchandise. Don’t miss out!" "附注 将免费赠送 Deriv 礼品。 千万不要错过！" $null $null $null
Replace-One "If you have questions, contact us " "如有任何疑问，请联系我们 " $null $null $null
Replace-One "[AFFILIATE MANAGER NO] (Viber/WhatsApp) | [AFFILIATE MANAGER EMAIL]" "[联盟经理编号] (Viber/WhatsApp) | [联盟经理电子邮件] (Viber/WhatsApp)" $null $null $null
Replace-One "SPANISH" "西班牙语" 1 $null $null
Replace-One "Back to" "转换回" $null $null $null
Replace-One "English" "英语" 1 13391121 1
